$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the params string for clf_fpipe_c: RandomForestClassifier -> ExtraTreesClassifier
$ws.Range("B2").Value = "{'clf': ExtraTreesClassifier(random_state=42), 'data_prep__numeric_pipe__data_missing__strategy': 'mean'}"

# Update the metric values to 1 (accuracy_test, f1_test, recall_test, roc_auc_test)
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 1
$ws.Range("F2").Value = 1
$ws.Range("H2").Value = 1
